$d = $word.ActiveDocument

# The bibliography paragraph currently holds all references concatenated
# into a single <w:t> run. We split it into separate <w:t> runs joined by
# manual line breaks (<w:br/>), one break inserted at each boundary
# between consecutive references.

$boundaries = @(
    @{ Find = " ed., Pearson, 2011.ZETTILI, N. Quantum "; Replace = " ed., Pearson, 2011.^lZETTILI, N. Quantum " },
    @{ Find = "ations, Wiley, 2009.CLAUDE COHEN-TANNOUDJI"; Replace = "ations, Wiley, 2009.^lCLAUDE COHEN-TANNOUDJI" },
    @{ Find = "iley and Sons, 1987.GASIOROWICZ, S., "; Replace = "iley and Sons, 1987.^lGASIOROWICZ, S., " },
    @{ Find = "bara Dois, RJ. 1979.FEYNMAN, R.P., "; Replace = "bara Dois, RJ. 1979.^lFEYNMAN, R.P., " },
    @{ Find = "ddison-Wesley, 1975.MERZBACHER, E., "; Replace = "ddison-Wesley, 1975.^lMERZBACHER, E., " },
    @{ Find = ", Nova Iorque, 1970.EISBERG, R.; RESNICK"; Replace = ", Nova Iorque, 1970.^lEISBERG, R.; RESNICK" }
)

foreach ($b in $boundaries) {
    $range = $d.Content
    $find = $range.Find
    $ok = $find.Execute($b.Find, $true, $false, $false, $false, $false, $true, 1, $false, $b.Replace, 2)
    if (-not $ok) {
        Write-Host "FAILED to find/replace:" $b.Find
    }
}

Write-Host "Done."
